# Apply updated Leve profit values across Sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6205.936
$ws.Range("I138").Value = 2237.8667
$ws.Range("J138").Value = 7150.7144
$ws.Range("K138").Value = 6713.6001
$ws.Range("L138").Value = 21452.1432
$ws.Range("M138").Value = -1573.6001
$ws.Range("N138").Value = -31732.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1130.1818
$ws.Range("I110").Value = 858.75
$ws.Range("K110").Value = 858.75
$ws.Range("M110").Value = 1186.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8785.579
$ws.Range("J20").Value = 7241.7
$ws.Range("L20").Value = 7241.7
$ws.Range("N20").Value = -7735.7
$ws.Range("H105").Value = 5337.6924
$ws.Range("I105").Value = 3969.2856
$ws.Range("J105").Value = 6934.1665
$ws.Range("K105").Value = 3969.2856
$ws.Range("L105").Value = 6934.1665
$ws.Range("M105").Value = -2222.2856
$ws.Range("N105").Value = -10428.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18521880
$ws.Range("I31").Value = 20411078
$ws.Range("J31").Value = 7734.8
$ws.Range("K31").Value = 20411078
$ws.Range("L31").Value = 7734.8
$ws.Range("M31").Value = -20410783
$ws.Range("N31").Value = -8324.799999999999
$ws.Range("H34").Value = 18521880
$ws.Range("I34").Value = 20411078
$ws.Range("J34").Value = 7734.8
$ws.Range("K34").Value = 20411078
$ws.Range("L34").Value = 7734.8
$ws.Range("M34").Value = -20410876
$ws.Range("N34").Value = -8138.8
$ws.Range("H58").Value = 3027.9443
$ws.Range("I58").Value = 2665.5
$ws.Range("K58").Value = 2665.5
$ws.Range("M58").Value = -2462.5
$ws.Range("H115").Value = 49832.668
$ws.Range("J115").Value = 49832.668
$ws.Range("L115").Value = 49832.668
$ws.Range("N115").Value = -52182.668
$ws.Range("H134").Value = 2150.3333
$ws.Range("I134").Value = 1920.8096
$ws.Range("K134").Value = 5762.4288
$ws.Range("M134").Value = -3227.4288
$ws.Range("H136").Value = 3027.9443
$ws.Range("I136").Value = 2665.5
$ws.Range("K136").Value = 7996.5
$ws.Range("M136").Value = -5446.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9216.556
$ws.Range("I3").Value = 3278.4285
$ws.Range("J3").Value = 30000
$ws.Range("K3").Value = 9835.2855
$ws.Range("L3").Value = 90000
$ws.Range("M3").Value = -9723.2855
$ws.Range("N3").Value = -90224
$ws.Range("H12").Value = 120.666664
$ws.Range("J12").Value = 109.14286
$ws.Range("L12").Value = 327.42858
$ws.Range("N12").Value = -673.42858
$ws.Range("H68").Value = 1397.6666
$ws.Range("I68").Value = 958.8
$ws.Range("K68").Value = 2876.4
$ws.Range("M68").Value = -2065.4
$ws.Range("H71").Value = 1397.6666
$ws.Range("I71").Value = 958.8
$ws.Range("K71").Value = 8629.199999999999
$ws.Range("M71").Value = -4573.199999999999
$ws.Range("H107").Value = 1942.7273
$ws.Range("J107").Value = 1637
$ws.Range("L107").Value = 4911
$ws.Range("N107").Value = -8751
$ws.Range("H122").Value = 1788.7778
$ws.Range("J122").Value = 1757
$ws.Range("L122").Value = 15813
$ws.Range("N122").Value = -20713
$ws.Range("H131").Value = 12822822
$ws.Range("J131").Value = 13891246
$ws.Range("L131").Value = 41673738
$ws.Range("N131").Value = -41683818
$ws.Range("H133").Value = 12903.692
$ws.Range("I133").Value = 7968.5
$ws.Range("J133").Value = 20800
$ws.Range("K133").Value = 23905.5
$ws.Range("L133").Value = 62400
$ws.Range("M133").Value = -18845.5
$ws.Range("N133").Value = -72520
$ws.Range("H140").Value = 4689.7144
$ws.Range("I140").Value = 4689.7144
$ws.Range("K140").Value = 14069.1432
$ws.Range("M140").Value = -8889.143199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7336.364
$ws.Range("I70").Value = 7212.5
$ws.Range("K70").Value = 7212.5
$ws.Range("M70").Value = -6942.5
$ws.Range("H73").Value = 7336.364
$ws.Range("I73").Value = 7212.5
$ws.Range("K73").Value = 7212.5
$ws.Range("M73").Value = -6276.5
$ws.Range("H102").Value = 10418722
$ws.Range("I102").Value = 12501917
$ws.Range("K102").Value = 12501917
$ws.Range("M102").Value = -12500295
$ws.Range("H107").Value = 903
$ws.Range("I107").Value = 705.1667
$ws.Range("K107").Value = 705.1667
$ws.Range("M107").Value = 1214.8333
$ws.Range("H113").Value = 2168.3076
$ws.Range("I113").Value = 1455.4445
$ws.Range("J113").Value = 3772.25
$ws.Range("K113").Value = 1455.4445
$ws.Range("L113").Value = 3772.25
$ws.Range("M113").Value = 714.5554999999999
$ws.Range("N113").Value = -8112.25
$ws.Range("H122").Value = 247901
$ws.Range("I122").Value = 402515.38
$ws.Range("J122").Value = 6316.0625
$ws.Range("K122").Value = 1207546.14
$ws.Range("L122").Value = 18948.1875
$ws.Range("M122").Value = -1205096.14
$ws.Range("N122").Value = -23848.1875
$ws.Range("H126").Value = 4512.4116
$ws.Range("I126").Value = 5737.909
$ws.Range("J126").Value = 2265.6667
$ws.Range("K126").Value = 17213.727
$ws.Range("L126").Value = 6797.000100000001
$ws.Range("M126").Value = -14743.727
$ws.Range("N126").Value = -11737.0001
$ws.Range("H132").Value = 102750.25
$ws.Range("I132").Value = 183993.19
$ws.Range("J132").Value = 3453.3333
$ws.Range("K132").Value = 551979.5700000001
$ws.Range("L132").Value = 10359.9999
$ws.Range("M132").Value = -549449.5700000001
$ws.Range("N132").Value = -15419.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 936.35
$ws.Range("I22").Value = 798.625
$ws.Range("J22").Value = 1487.25
$ws.Range("K22").Value = 798.625
$ws.Range("L22").Value = 1487.25
$ws.Range("M22").Value = -503.625
$ws.Range("N22").Value = -2077.25
$ws.Range("H27").Value = 936.35
$ws.Range("I27").Value = 798.625
$ws.Range("J27").Value = 1487.25
$ws.Range("K27").Value = 798.625
$ws.Range("L27").Value = 1487.25
$ws.Range("M27").Value = -691.625
$ws.Range("N27").Value = -1701.25
$ws.Range("H33").Value = 15000000
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H47").Value = 38435.43
$ws.Range("J47").Value = 41508.168
$ws.Range("L47").Value = 41508.168
$ws.Range("N47").Value = -42488.168
$ws.Range("H52").Value = 38435.43
$ws.Range("J52").Value = 41508.168
$ws.Range("L52").Value = 41508.168
$ws.Range("N52").Value = -41974.168
$ws.Range("H55").Value = 610.5294
$ws.Range("J55").Value = 1012.3333
$ws.Range("L55").Value = 1012.3333
$ws.Range("N55").Value = -1358.3333
$ws.Range("H76").Value = 48341.75
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 48341.75
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 48341.75
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -49017.75
$ws.Range("H79").Value = 48341.75
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 48341.75
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 48341.75
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -50681.75
$ws.Range("H132").Value = 2672.85
$ws.Range("I132").Value = 2620.1912
$ws.Range("K132").Value = 7860.573600000001
$ws.Range("M132").Value = -5330.573600000001
$ws.Range("H140").Value = 67247.46000000001
$ws.Range("J140").Value = 67247.46000000001
$ws.Range("L140").Value = 67247.46000000001
$ws.Range("N140").Value = -77607.46000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9639.467000000001
$ws.Range("I81").Value = 8699.429
$ws.Range("J81").Value = 10462
$ws.Range("K81").Value = 17398.858
$ws.Range("L81").Value = 20924
$ws.Range("M81").Value = -16337.858
$ws.Range("N81").Value = -23046
$ws.Range("H84").Value = 9639.467000000001
$ws.Range("I84").Value = 8699.429
$ws.Range("J84").Value = 10462
$ws.Range("K84").Value = 86994.29000000001
$ws.Range("L84").Value = 104620
$ws.Range("M84").Value = -81690.29000000001
$ws.Range("N84").Value = -115228
$ws.Range("H113").Value = 750.8333
$ws.Range("I113").Value = 543
$ws.Range("K113").Value = 1629
$ws.Range("M113").Value = 541
$ws.Range("H126").Value = 52634190
$ws.Range("I126").Value = 58825980
$ws.Range("J126").Value = 3998
$ws.Range("K126").Value = 176477940
$ws.Range("L126").Value = 11994
$ws.Range("M126").Value = -176475470
$ws.Range("N126").Value = -16934
